# Automated task-scheduler style update: refresh the last existing reading's
# timestamp precision and append the new sensor reading as the next row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: floating-point precision refresh of the stored timestamp (same
# moment, 11:00:16, just re-serialized with slightly different precision).
$ws.Range("A11").Value = 45865.45852378472

# Row 12: new sensor reading appended by the scheduled task.
$ws.Range("A12").Value = 45865.54202206938
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat
$ws.Range("B12").Value = 2025
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 16.5
$ws.Range("E12").Value = 81.81999999999999
$ws.Range("F12").Value = 638.09
$ws.Range("G12").Value = 13.4
$ws.Range("H12").Value = "SE"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "13:00:30"
